# Add two new review rows (40-41) plus three blank-but-formatted rows (42-44)
# to Sheet1, matching the appended reviews for com.hamxa.shaynachim / bitcoin,
# and move the selection to F42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40 -----------------------------------------------------------
# Seed formatting by copying the last existing data row (39), then
# overwrite the values. This reproduces the per-column styles
# (A/E/F/G = style 0 or 1, C/D = style 2) without hard-coding style ids.
$ws.Range("A39:G39").Copy()
$ws.Range("A40:G40").PasteSpecial(-4122)

$ws.Range("A40").Value = "com.hamxa.shaynachim"
$ws.Range("B40").Value = "bitcoin"
$ws.Range("C40").Value = "sixsevensix67676@gmail.com"
$ws.Range("D40").Value = "stevewonder3001@gmail.com"
$ws.Range("E40").Value = "27/5/2019 15:59"
$ws.Range("F40").Value = "beginner guide that is worth 5 star for sure. Thanks"
$ws.Range("G40").Value = "yes"

$ws.Hyperlinks.Add($ws.Range("C40"), "mailto:sixsevensix67676@gmail.com", "", "", "sixsevensix67676@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D40"), "mailto:stevewonder3001@gmail.com", "", "", "stevewonder3001@gmail.com")

# Adding a hyperlink re-styles the cell with the built-in "Hyperlink" style;
# restore the sheet's normal wrapped-text style (same as every other
# email/recovery cell, hyperlinked or not).
$ws.Range("C39:D39").Copy()
$ws.Range("C40:D40").PasteSpecial(-4122)

# --- Row 41 -----------------------------------------------------------
$ws.Range("A39:G39").Copy()
$ws.Range("A41:G41").PasteSpecial(-4122)

$ws.Range("A41").Value = "com.hamxa.shaynachim"
$ws.Range("B41").Value = "bitcoin"
$ws.Range("C41").Value = "dony1098765432@gmail.com"
$ws.Range("D41").Value = "sixsevensix67676@gmail.com"
$ws.Range("E41").Value = "27/5/2019 15:59"
$ws.Range("F41").Value = "welcome to the best guide about bitcoin this year"
$ws.Range("G41").Value = "yes"

$ws.Hyperlinks.Add($ws.Range("C41"), "mailto:dony1098765432@gmail.com", "", "", "dony1098765432@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D41"), "mailto:sixsevensix67676@gmail.com", "", "", "sixsevensix67676@gmail.com")

# Restore styling clobbered by Hyperlinks.Add (see row 40 above).
$ws.Range("C39:D39").Copy()
$ws.Range("C41:D41").PasteSpecial(-4122)

# --- Rows 42-44: blank trailer rows, formatted like the C/D data cells --
$ws.Range("C39:D39").Copy()
$ws.Range("C42:D42").PasteSpecial(-4122)
$ws.Range("C39:D39").Copy()
$ws.Range("C43:D43").PasteSpecial(-4122)
$ws.Range("C39:D39").Copy()
$ws.Range("C44:D44").PasteSpecial(-4122)

$ws.Rows.Item(42).RowHeight = 13.8
$ws.Rows.Item(43).RowHeight = 13.8
$ws.Rows.Item(44).RowHeight = 13.8

# Move the visible selection, like in the saved workbook.
[void]$ws.Range("F42").Select()

Write-Host "done"
